$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 14 ("Nov 25 - Nov 27") currently lists topic "Lambda functions".
# Replace it with "Code Smells" per the commit "added code smell to topics".
$ws.Range("C14").Value = "Code Smells"
